# Trade #299 closed at 2026-02-18 01:43:54 - unknown UNKNOWN +0.000%
#
# Updates the "Summary", "Strategy Status" and "All Trades" sheets to
# reflect a closed MarketMaking trade, and appends the newly-opened trades
# that were recorded afterwards (both to "All Trades" and to each
# strategy's own sheet).

$wb = $excel.ActiveWorkbook

function Set-DateCell($ws, $row, $col, $text) {
    # Excel's COM layer auto-parses strings that look like dates into date
    # serials. Forcing the cell to Text format before assigning the value
    # keeps it as the literal string, matching the source data.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.15
$summary.Range("B4").Value = 0.25
$summary.Range("B5").Value = 0.02
$summary.Range("B6").Value = 327
$summary.Range("B7").Value = 123
$summary.Range("B9").Value = 37.61

# ---------------------------------------------------------------------
# Strategy Status sheet (row 6 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.08
$status.Range("D6").Value = 194
$status.Range("E6").Value = -0.73
$status.Range("F6").Value = -0.92
$status.Range("G6").Value = 34.54

# ---------------------------------------------------------------------
# All Trades sheet - trade #327 (row 328) closes out
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(328, 7).Value = 0.99          # G: Exit Price
$allTrades.Cells.Item(328, 8).Value = "CLOSED"      # H: Status
$allTrades.Cells.Item(328, 9).Value = 2.0619        # I: P&L %
$allTrades.Cells.Item(328, 10).Value = 0.02         # J: P&L $
$allTrades.Cells.Item(328, 11).Value = 99.08        # K: Capital After
$allTrades.Cells.Item(328, 12).Value = "early_exit" # L: Exit Reason
$allTrades.Cells.Item(328, 13).Value = 0.18         # M: Duration (min)

# ---------------------------------------------------------------------
# All Trades sheet - 4 newly opened trades appended (rows 357-360)
# ---------------------------------------------------------------------
function Add-AllTradesRow($row, $tradeNum, $date, $time, $strategy, $side, $entry, $status, $plPct, $plUsd, $capAfter, $exitReason, $duration, $entrySlip, $exitSlip, $confidence, $entryReason) {
    $allTrades.Cells.Item($row, 1).Value = $tradeNum
    Set-DateCell $allTrades $row 2 $date
    $allTrades.Cells.Item($row, 3).Value = $time
    $allTrades.Cells.Item($row, 4).Value = $strategy
    $allTrades.Cells.Item($row, 5).Value = $side
    $allTrades.Cells.Item($row, 6).Value = $entry
    # G (Exit Price) left blank - trade is still OPEN
    $allTrades.Cells.Item($row, 8).Value = $status
    $allTrades.Cells.Item($row, 9).Value = $plPct
    $allTrades.Cells.Item($row, 10).Value = $plUsd
    $allTrades.Cells.Item($row, 11).Value = $capAfter
    # L (Exit Reason) left blank - trade is still OPEN
    $allTrades.Cells.Item($row, 13).Value = $duration
    $allTrades.Cells.Item($row, 14).Value = $entrySlip
    $allTrades.Cells.Item($row, 15).Value = $exitSlip
    $allTrades.Cells.Item($row, 16).Value = $confidence
    $allTrades.Cells.Item($row, 17).Value = $entryReason
}

Add-AllTradesRow 357 356 "2026-02-18" "01:43:47" "momentum" "DOWN" 0.97 "OPEN" 0 0 99.67024513670182 "" 0 0 0 0.9 "Downward momentum: -12.621% over 10 samples"
Add-AllTradesRow 358 357 "2026-02-18" "01:43:47" "HighProbConvergence" "UP" 0.03 "OPEN" 0 0 100.1259031022653 "" 0 0 0 0.95 "Mean reversion UP: price 10.94% below mean (z=-2.38)"
Add-AllTradesRow 359 358 "2026-02-18" "01:43:48" "MarketMaking" "DOWN" 0.97 "OPEN" 0 0 99.06245085365964 "" 0 0 0 0.65 "Wide spread capture: 444 bps vs avg 209 bps"
Add-AllTradesRow 360 359 "2026-02-18" "01:43:49" "EMAArbitrage" "DOWN" 0.98 "OPEN" 0 0 100.270616878256 "" 0 0 0 0.6262 "EMA:down, RSI:0.0, ROC:-12.62% | 2/3 DOWN"

# ---------------------------------------------------------------------
# momentum sheet - trade #356 appended (row 67)
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L%, J P&L$, K Capital After,
#          L Entry Slippage, M Exit Slippage, N Confidence,
#          O Entry Reason, P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(67, 1).Value = 356
Set-DateCell $momentum 67 2 "2026-02-18"
$momentum.Cells.Item(67, 3).Value = "01:43:47"
$momentum.Cells.Item(67, 4).Value = "momentum"
$momentum.Cells.Item(67, 5).Value = "DOWN"
$momentum.Cells.Item(67, 6).Value = 0.97
# G (Exit Price) left blank - trade is still OPEN
$momentum.Cells.Item(67, 8).Value = "OPEN"
$momentum.Cells.Item(67, 9).Value = 0
$momentum.Cells.Item(67, 10).Value = 0
$momentum.Cells.Item(67, 11).Value = 99.67024513670182
$momentum.Cells.Item(67, 12).Value = 0
$momentum.Cells.Item(67, 13).Value = 0
$momentum.Cells.Item(67, 14).Value = 0.9
$momentum.Cells.Item(67, 15).Value = "Downward momentum: -12.621% over 10 samples"
# P (Exit Reason) left blank - trade is still OPEN
$momentum.Cells.Item(67, 17).Value = 0

# ---------------------------------------------------------------------
# HighProbConvergence sheet - trade #357 appended (row 32)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(32, 1).Value = 357
Set-DateCell $hpc 32 2 "2026-02-18"
$hpc.Cells.Item(32, 3).Value = "01:43:47"
$hpc.Cells.Item(32, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(32, 5).Value = "UP"
$hpc.Cells.Item(32, 6).Value = 0.03
# G (Exit Price) left blank - trade is still OPEN
$hpc.Cells.Item(32, 8).Value = "OPEN"
$hpc.Cells.Item(32, 9).Value = 0
$hpc.Cells.Item(32, 10).Value = 0
$hpc.Cells.Item(32, 11).Value = 100.1259031022653
$hpc.Cells.Item(32, 12).Value = 0
$hpc.Cells.Item(32, 13).Value = 0
$hpc.Cells.Item(32, 14).Value = 0.95
$hpc.Cells.Item(32, 15).Value = "Mean reversion UP: price 10.94% below mean (z=-2.38)"
# P (Exit Reason) left blank - trade is still OPEN
$hpc.Cells.Item(32, 17).Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - trade #327 (row 195) closes out, matching
# "All Trades" row 328, plus trade #358 appended (row 221)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(195, 7).Value = 0.99          # G: Exit Price
$mm.Cells.Item(195, 8).Value = "CLOSED"      # H: Status
$mm.Cells.Item(195, 9).Value = 2.0619        # I: P&L %
$mm.Cells.Item(195, 10).Value = 0.02         # J: P&L $
$mm.Cells.Item(195, 11).Value = 99.08        # K: Capital After
$mm.Cells.Item(195, 16).Value = "early_exit" # P: Exit Reason
$mm.Cells.Item(195, 17).Value = 0.18         # Q: Duration (min)

$mm.Cells.Item(221, 1).Value = 358
Set-DateCell $mm 221 2 "2026-02-18"
$mm.Cells.Item(221, 3).Value = "01:43:48"
$mm.Cells.Item(221, 4).Value = "MarketMaking"
$mm.Cells.Item(221, 5).Value = "DOWN"
$mm.Cells.Item(221, 6).Value = 0.97
# G (Exit Price) left blank - trade is still OPEN
$mm.Cells.Item(221, 8).Value = "OPEN"
$mm.Cells.Item(221, 9).Value = 0
$mm.Cells.Item(221, 10).Value = 0
$mm.Cells.Item(221, 11).Value = 99.06245085365964
$mm.Cells.Item(221, 12).Value = 0
$mm.Cells.Item(221, 13).Value = 0
$mm.Cells.Item(221, 14).Value = 0.65
$mm.Cells.Item(221, 15).Value = "Wide spread capture: 444 bps vs avg 209 bps"
# P (Exit Reason) left blank - trade is still OPEN
$mm.Cells.Item(221, 17).Value = 0

# ---------------------------------------------------------------------
# EMAArbitrage sheet - trade #359 appended (row 14)
# ---------------------------------------------------------------------
$ema = $wb.Worksheets.Item("EMAArbitrage")
$ema.Cells.Item(14, 1).Value = 359
Set-DateCell $ema 14 2 "2026-02-18"
$ema.Cells.Item(14, 3).Value = "01:43:49"
$ema.Cells.Item(14, 4).Value = "EMAArbitrage"
$ema.Cells.Item(14, 5).Value = "DOWN"
$ema.Cells.Item(14, 6).Value = 0.98
# G (Exit Price) left blank - trade is still OPEN
$ema.Cells.Item(14, 8).Value = "OPEN"
$ema.Cells.Item(14, 9).Value = 0
$ema.Cells.Item(14, 10).Value = 0
$ema.Cells.Item(14, 11).Value = 100.270616878256
$ema.Cells.Item(14, 12).Value = 0
$ema.Cells.Item(14, 13).Value = 0
$ema.Cells.Item(14, 14).Value = 0.6262
$ema.Cells.Item(14, 15).Value = "EMA:down, RSI:0.0, ROC:-12.62% | 2/3 DOWN"
# P (Exit Reason) left blank - trade is still OPEN
$ema.Cells.Item(14, 17).Value = 0
